# ----------------------------------------------------------------------------
# Fix Training Data Issue
# The box-score stats in this sheet were pulled one calendar day off because of
# how NBA.com displayed game dates; this corrects the affected numeric stats/
# ranks for every team row and normalizes the Date column (BF) from the old
# "1-18-2008-09" label to the ISO-style "2009-01-18".
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Date column (BF) holds text that looks like a date (e.g. "2009-01-18").
# Force the cells to text format first so Excel stores the literal string
# instead of silently converting it to a date serial number.
$ws.Range("BF2").NumberFormat = "@"
$ws.Range("BF3").NumberFormat = "@"
$ws.Range("BF4").NumberFormat = "@"
$ws.Range("BF5").NumberFormat = "@"
$ws.Range("BF6").NumberFormat = "@"
$ws.Range("BF7").NumberFormat = "@"
$ws.Range("BF8").NumberFormat = "@"
$ws.Range("BF9").NumberFormat = "@"
$ws.Range("BF10").NumberFormat = "@"
$ws.Range("BF11").NumberFormat = "@"
$ws.Range("BF12").NumberFormat = "@"
$ws.Range("BF13").NumberFormat = "@"
$ws.Range("BF14").NumberFormat = "@"
$ws.Range("BF15").NumberFormat = "@"
$ws.Range("BF16").NumberFormat = "@"
$ws.Range("BF17").NumberFormat = "@"
$ws.Range("BF18").NumberFormat = "@"
$ws.Range("BF19").NumberFormat = "@"
$ws.Range("BF20").NumberFormat = "@"
$ws.Range("BF21").NumberFormat = "@"
$ws.Range("BF22").NumberFormat = "@"
$ws.Range("BF23").NumberFormat = "@"
$ws.Range("BF24").NumberFormat = "@"
$ws.Range("BF25").NumberFormat = "@"
$ws.Range("BF26").NumberFormat = "@"
$ws.Range("BF27").NumberFormat = "@"
$ws.Range("BF28").NumberFormat = "@"
$ws.Range("BF29").NumberFormat = "@"
$ws.Range("BF30").NumberFormat = "@"
$ws.Range("BF31").NumberFormat = "@"

# Row 2
$ws.Range("AF2").Value = 9
$ws.Range("AI2").Value = 21
$ws.Range("AK2").Value = 10
$ws.Range("AN2").Value = 9
$ws.Range("AP2").Value = 16
$ws.Range("AR2").Value = 18
$ws.Range("AU2").Value = 9
$ws.Range("AX2").Value = 18
$ws.Range("BB2").Value = 15
$ws.Range("BF2").Value = "2009-01-18"

# Row 3
$ws.Range("D3").Value = 42
$ws.Range("E3").Value = 33
$ws.Range("G3").Value = 0.786
$ws.Range("I3").Value = 36.6
$ws.Range("J3").Value = 76.3
$ws.Range("K3").Value = 0.48
$ws.Range("L3").Value = 6.4
$ws.Range("N3").Value = 0.374
$ws.Range("P3").Value = 27.5
$ws.Range("Q3").Value = 0.771
$ws.Range("S3").Value = 32
$ws.Range("T3").Value = 42.8
$ws.Range("U3").Value = 22.3
$ws.Range("V3").Value = 16.2
$ws.Range("X3").Value = 4.9
$ws.Range("Y3").Value = 4.6
$ws.Range("Z3").Value = 23.2
$ws.Range("AB3").Value = 100.9
$ws.Range("AC3").Value = 9.199999999999999
$ws.Range("AD3").Value = 2
$ws.Range("AE3").Value = 1
$ws.Range("AH3").Value = 13
$ws.Range("AN3").Value = 12
$ws.Range("AQ3").Value = 12
$ws.Range("AR3").Value = 17
$ws.Range("AT3").Value = 5
$ws.Range("AX3").Value = 16
$ws.Range("AY3").Value = 12
$ws.Range("BF3").Value = "2009-01-18"

# Row 4
$ws.Range("AD4").Value = 13
$ws.Range("AK4").Value = 21
$ws.Range("BF4").Value = "2009-01-18"

# Row 5
$ws.Range("D5").Value = 41
$ws.Range("E5").Value = 18
$ws.Range("G5").Value = 0.439
$ws.Range("I5").Value = 37.1
$ws.Range("J5").Value = 83.8
$ws.Range("K5").Value = 0.443
$ws.Range("L5").Value = 6.2
$ws.Range("M5").Value = 16.4
$ws.Range("N5").Value = 0.38
$ws.Range("O5").Value = 18.7
$ws.Range("P5").Value = 23.5
$ws.Range("Q5").Value = 0.798
$ws.Range("R5").Value = 11.6
$ws.Range("S5").Value = 30.1
$ws.Range("T5").Value = 41.7
$ws.Range("U5").Value = 20.5
$ws.Range("V5").Value = 14.9
$ws.Range("Y5").Value = 5.6
$ws.Range("Z5").Value = 22.1
$ws.Range("AA5").Value = 20.4
$ws.Range("AB5").Value = 99.09999999999999
$ws.Range("AC5").Value = -3.1
$ws.Range("AD5").Value = 3
$ws.Range("AF5").Value = 18
$ws.Range("AH5").Value = 8
$ws.Range("AN5").Value = 8
$ws.Range("AO5").Value = 17
$ws.Range("AP5").Value = 21
$ws.Range("AS5").Value = 15
$ws.Range("AW5").Value = 12
$ws.Range("AX5").Value = 6
$ws.Range("BA5").Value = 22
$ws.Range("BF5").Value = "2009-01-18"

# Row 6
$ws.Range("D6").Value = 38
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 0.8159999999999999
$ws.Range("I6").Value = 37.2
$ws.Range("J6").Value = 77.8
$ws.Range("K6").Value = 0.478
$ws.Range("L6").Value = 7.4
$ws.Range("M6").Value = 20.5
$ws.Range("N6").Value = 0.36
$ws.Range("O6").Value = 19.2
$ws.Range("P6").Value = 25.4
$ws.Range("Q6").Value = 0.758
$ws.Range("R6").Value = 10.5
$ws.Range("U6").Value = 20.2
$ws.Range("V6").Value = 13.3
$ws.Range("X6").Value = 5.8
$ws.Range("Y6").Value = 3.6
$ws.Range("Z6").Value = 20.7
$ws.Range("AA6").Value = 21.1
$ws.Range("AB6").Value = 100.9
$ws.Range("AC6").Value = 11.9
$ws.Range("AD6").Value = 27
$ws.Range("AL6").Value = 8
$ws.Range("AM6").Value = 7
$ws.Range("AO6").Value = 13
$ws.Range("AP6").Value = 10
$ws.Range("AQ6").Value = 20
$ws.Range("AU6").Value = 19
$ws.Range("AZ6").Value = 12
$ws.Range("BA6").Value = 16
$ws.Range("BF6").Value = "2009-01-18"

# Row 7
$ws.Range("AD7").Value = 13
$ws.Range("AH7").Value = 16
$ws.Range("AO7").Value = 27
$ws.Range("AU7").Value = 8
$ws.Range("BF7").Value = "2009-01-18"

# Row 8
$ws.Range("AD8").Value = 3
$ws.Range("AJ8").Value = 23
$ws.Range("AN8").Value = 10
$ws.Range("AT8").Value = 18
$ws.Range("BA8").Value = 4
$ws.Range("BC8").Value = 7
$ws.Range("BF8").Value = "2009-01-18"

# Row 9
$ws.Range("AE9").Value = 13
$ws.Range("AR9").Value = 18
$ws.Range("BA9").Value = 23
$ws.Range("BC9").Value = 16
$ws.Range("BF9").Value = "2009-01-18"

# Row 10
$ws.Range("AD10").Value = 3
$ws.Range("AU10").Value = 15
$ws.Range("AW10").Value = 6
$ws.Range("BF10").Value = "2009-01-18"

# Row 11
$ws.Range("AD11").Value = 3
$ws.Range("AF11").Value = 9
$ws.Range("AG11").Value = 8
$ws.Range("AO11").Value = 9
$ws.Range("BA11").Value = 17
$ws.Range("BC11").Value = 9
$ws.Range("BF11").Value = "2009-01-18"

# Row 12
$ws.Range("AD12").Value = 13
$ws.Range("AL12").Value = 7
$ws.Range("AO12").Value = 16
$ws.Range("AP12").Value = 22
$ws.Range("AU12").Value = 2
$ws.Range("AW12").Value = 16
$ws.Range("AZ12").Value = 29
$ws.Range("BA12").Value = 12
$ws.Range("BF12").Value = "2009-01-18"

# Row 13
$ws.Range("AP13").Value = 24
$ws.Range("BF13").Value = "2009-01-18"

# Row 14
$ws.Range("D14").Value = 39
$ws.Range("E14").Value = 31
$ws.Range("G14").Value = 0.795
$ws.Range("J14").Value = 83.7
$ws.Range("M14").Value = 18.4
$ws.Range("N14").Value = 0.383
$ws.Range("O14").Value = 21.3
$ws.Range("P14").Value = 27.6
$ws.Range("Q14").Value = 0.771
$ws.Range("R14").Value = 12.1
$ws.Range("S14").Value = 31.9
$ws.Range("T14").Value = 44
$ws.Range("U14").Value = 23.1
$ws.Range("X14").Value = 5.1
$ws.Range("Z14").Value = 20.2
$ws.Range("AA14").Value = 23.1
$ws.Range("AB14").Value = 107.6
$ws.Range("AC14").Value = 8.199999999999999
$ws.Range("AD14").Value = 18
$ws.Range("AE14").Value = 3
$ws.Range("AN14").Value = 6
$ws.Range("AO14").Value = 4
$ws.Range("AQ14").Value = 13
$ws.Range("AU14").Value = 3
$ws.Range("AX14").Value = 12
$ws.Range("BF14").Value = "2009-01-18"

# Row 15
$ws.Range("BF15").Value = "2009-01-18"

# Row 16
$ws.Range("D16").Value = 39
$ws.Range("E16").Value = 21
$ws.Range("G16").Value = 0.538
$ws.Range("I16").Value = 36.2
$ws.Range("J16").Value = 81.09999999999999
$ws.Range("K16").Value = 0.447
$ws.Range("L16").Value = 6.8
$ws.Range("M16").Value = 19.2
$ws.Range("N16").Value = 0.355
$ws.Range("O16").Value = 17.2
$ws.Range("P16").Value = 23.3
$ws.Range("Q16").Value = 0.738
$ws.Range("R16").Value = 10.8
$ws.Range("S16").Value = 29.5
$ws.Range("T16").Value = 40.3
$ws.Range("U16").Value = 19.9
$ws.Range("X16").Value = 5.5
$ws.Range("Y16").Value = 4.3
$ws.Range("AA16").Value = 19.8
$ws.Range("AB16").Value = 96.5
$ws.Range("AC16").Value = -0.2
$ws.Range("AD16").Value = 18
$ws.Range("AE16").Value = 15
$ws.Range("AH16").Value = 14
$ws.Range("AI16").Value = 19
$ws.Range("AK16").Value = 20
$ws.Range("AL16").Value = 14
$ws.Range("AN16").Value = 18
$ws.Range("AO16").Value = 26
$ws.Range("AR16").Value = 16
$ws.Range("AS16").Value = 23
$ws.Range("AU16").Value = 25
$ws.Range("AW16").Value = 5
$ws.Range("AX16").Value = 7
$ws.Range("BC16").Value = 17
$ws.Range("BF16").Value = "2009-01-18"

# Row 17
$ws.Range("D17").Value = 43
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 0.465
$ws.Range("H17").Value = 48.3
$ws.Range("J17").Value = 81.7
$ws.Range("K17").Value = 0.447
$ws.Range("N17").Value = 0.349
$ws.Range("O17").Value = 19.5
$ws.Range("P17").Value = 25
$ws.Range("Q17").Value = 0.782
$ws.Range("S17").Value = 29.3
$ws.Range("T17").Value = 41.8
$ws.Range("V17").Value = 14.6
$ws.Range("Z17").Value = 24
$ws.Range("AA17").Value = 22.6
$ws.Range("AB17").Value = 98
$ws.Range("AC17").Value = 0.5
$ws.Range("AF17").Value = 18
$ws.Range("AH17").Value = 19
$ws.Range("AJ17").Value = 10
$ws.Range("AK17").Value = 19
$ws.Range("AO17").Value = 10
$ws.Range("AS17").Value = 26
$ws.Range("AW17").Value = 17
$ws.Range("AZ17").Value = 30
$ws.Range("BB17").Value = 16
$ws.Range("BC17").Value = 15
$ws.Range("BF17").Value = "2009-01-18"

# Row 18
$ws.Range("AD18").Value = 27
$ws.Range("AF18").Value = 24
$ws.Range("AP18").Value = 11
$ws.Range("AQ18").Value = 19
$ws.Range("AU18").Value = 18
$ws.Range("BF18").Value = "2009-01-18"

# Row 19
$ws.Range("AD19").Value = 3
$ws.Range("AM19").Value = 6
$ws.Range("AN19").Value = 11
$ws.Range("BF19").Value = "2009-01-18"

# Row 20
$ws.Range("AK20").Value = 12
$ws.Range("AO20").Value = 21
$ws.Range("AS20").Value = 25
$ws.Range("BF20").Value = "2009-01-18"

# Row 21
$ws.Range("AG21").Value = 22
$ws.Range("AN21").Value = 17
$ws.Range("AU21").Value = 9
$ws.Range("AY21").Value = 22
$ws.Range("BF21").Value = "2009-01-18"

# Row 22
$ws.Range("D22").Value = 41
$ws.Range("F22").Value = 33
$ws.Range("G22").Value = 0.195
$ws.Range("I22").Value = 36.4
$ws.Range("J22").Value = 81.7
$ws.Range("K22").Value = 0.445
$ws.Range("N22").Value = 0.368
$ws.Range("O22").Value = 18.7
$ws.Range("P22").Value = 24.6
$ws.Range("Q22").Value = 0.763
$ws.Range("R22").Value = 11.7
$ws.Range("T22").Value = 42.7
$ws.Range("U22").Value = 20.1
$ws.Range("X22").Value = 4.6
$ws.Range("Y22").Value = 5.3
$ws.Range("Z22").Value = 20.9
$ws.Range("AB22").Value = 95.40000000000001
$ws.Range("AC22").Value = -6.9
$ws.Range("AD22").Value = 3
$ws.Range("AH22").Value = 21
$ws.Range("AI22").Value = 16
$ws.Range("AJ22").Value = 9
$ws.Range("AK22").Value = 23
$ws.Range("AP22").Value = 17
$ws.Range("AQ22").Value = 15
$ws.Range("AU22").Value = 22
$ws.Range("AY22").Value = 21
$ws.Range("AZ22").Value = 15
$ws.Range("BA22").Value = 24
$ws.Range("BF22").Value = "2009-01-18"

# Row 23
$ws.Range("AD23").Value = 3
$ws.Range("AJ23").Value = 22
$ws.Range("AO23").Value = 14
$ws.Range("AT23").Value = 4
$ws.Range("AY23").Value = 3
$ws.Range("BF23").Value = "2009-01-18"

# Row 24
$ws.Range("AD24").Value = 13
$ws.Range("AK24").Value = 9
$ws.Range("AO24").Value = 22
$ws.Range("BC24").Value = 14
$ws.Range("BF24").Value = "2009-01-18"

# Row 25
$ws.Range("D25").Value = 37
$ws.Range("E25").Value = 22
$ws.Range("G25").Value = 0.595
$ws.Range("I25").Value = 38.6
$ws.Range("K25").Value = 0.498
$ws.Range("L25").Value = 6.9
$ws.Range("N25").Value = 0.388
$ws.Range("O25").Value = 20.4
$ws.Range("P25").Value = 26.9
$ws.Range("Q25").Value = 0.76
$ws.Range("R25").Value = 9.6
$ws.Range("S25").Value = 31.5
$ws.Range("T25").Value = 41.1
$ws.Range("U25").Value = 21.5
$ws.Range("V25").Value = 16.2
$ws.Range("W25").Value = 6.4
$ws.Range("X25").Value = 4.9
$ws.Range("Z25").Value = 20.5
$ws.Range("AA25").Value = 22.1
$ws.Range("AB25").Value = 104.5
$ws.Range("AD25").Value = 29
$ws.Range("AE25").Value = 13
$ws.Range("AL25").Value = 13
$ws.Range("AQ25").Value = 17
$ws.Range("AS25").Value = 7
$ws.Range("AT25").Value = 17
$ws.Range("AU25").Value = 11
$ws.Range("AX25").Value = 15
$ws.Range("BF25").Value = "2009-01-18"

# Row 26
$ws.Range("D26").Value = 40
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 0.6
$ws.Range("I26").Value = 35.9
$ws.Range("J26").Value = 78.8
$ws.Range("K26").Value = 0.456
$ws.Range("L26").Value = 7.6
$ws.Range("M26").Value = 19.7
$ws.Range("N26").Value = 0.383
$ws.Range("O26").Value = 18.6
$ws.Range("P26").Value = 24.2
$ws.Range("Q26").Value = 0.767
$ws.Range("S26").Value = 27.7
$ws.Range("T26").Value = 40.7
$ws.Range("V26").Value = 12.9
$ws.Range("W26").Value = 6.8
$ws.Range("X26").Value = 5
$ws.Range("Y26").Value = 3.8
$ws.Range("Z26").Value = 20.8
$ws.Range("AA26").Value = 21.3
$ws.Range("AB26").Value = 97.90000000000001
$ws.Range("AC26").Value = 2.6
$ws.Range("AD26").Value = 13
$ws.Range("AF26").Value = 9
$ws.Range("AG26").Value = 9
$ws.Range("AH26").Value = 11
$ws.Range("AI26").Value = 22
$ws.Range("AJ26").Value = 21
$ws.Range("AK26").Value = 13
$ws.Range("AN26").Value = 5
$ws.Range("AO26").Value = 20
$ws.Range("AQ26").Value = 14
$ws.Range("AU26").Value = 16
$ws.Range("AX26").Value = 13
$ws.Range("AY26").Value = 4
$ws.Range("AZ26").Value = 14
$ws.Range("BA26").Value = 13
$ws.Range("BB26").Value = 17
$ws.Range("BC26").Value = 10
$ws.Range("BF26").Value = "2009-01-18"

# Row 27
$ws.Range("AD27").Value = 3
$ws.Range("AI27").Value = 18
$ws.Range("AK27").Value = 22
$ws.Range("AS27").Value = 24
$ws.Range("AU27").Value = 23
$ws.Range("BF27").Value = "2009-01-18"

# Row 28
$ws.Range("AQ28").Value = 16
$ws.Range("AS28").Value = 8
$ws.Range("BC28").Value = 8
$ws.Range("BF28").Value = "2009-01-18"

# Row 29
$ws.Range("D29").Value = 41
$ws.Range("F29").Value = 25
$ws.Range("G29").Value = 0.39
$ws.Range("I29").Value = 35.8
$ws.Range("K29").Value = 0.457
$ws.Range("N29").Value = 0.38
$ws.Range("Q29").Value = 0.826
$ws.Range("R29").Value = 8.9
$ws.Range("S29").Value = 30.3
$ws.Range("T29").Value = 39.2
$ws.Range("W29").Value = 6.5
$ws.Range("Y29").Value = 4.6
$ws.Range("Z29").Value = 19.5
$ws.Range("AB29").Value = 97.5
$ws.Range("AC29").Value = -2.4
$ws.Range("AD29").Value = 3
$ws.Range("AF29").Value = 22
$ws.Range("AG29").Value = 21
$ws.Range("AH29").Value = 17
$ws.Range("AI29").Value = 23
$ws.Range("AJ29").Value = 24
$ws.Range("AK29").Value = 11
$ws.Range("AN29").Value = 7
$ws.Range("AS29").Value = 13
$ws.Range("AU29").Value = 7
$ws.Range("AX29").Value = 14
$ws.Range("AY29").Value = 13
$ws.Range("BF29").Value = "2009-01-18"

# Row 30
$ws.Range("AD30").Value = 3
$ws.Range("AO30").Value = 3
$ws.Range("AQ30").Value = 11
$ws.Range("AS30").Value = 22
$ws.Range("AX30").Value = 17
$ws.Range("BC30").Value = 6
$ws.Range("BF30").Value = "2009-01-18"

# Row 31
$ws.Range("AU31").Value = 17
$ws.Range("AZ31").Value = 11
$ws.Range("BF31").Value = "2009-01-18"
